# Applies the "fix vacc data and small improvements for report" commit.
#
# Summary of changes:
#  1) Sheet "Todesfälle und Fallsterblichkei" (index 5): fix sign/space typos
#     in the "Veraenderung" (D) column for the Uebersterblichkeit block.
#  2) Sheet "Geimpfte Personen" (index 9): refreshed to 4.5. -> 5.5. data,
#     gains a new "Vorwoche" column (B) + renamed/re-added change column (D).
#  3) Sheet "Impffortschritt" (index 10): new "Veraenderung" column (D) with
#     percentage deltas.
#  4) Sheet "Regional Geimpfte" (index 11): refreshed "Gesamt min. 1x" (C) and
#     "Gesamt vollst." (D) percentages for all 17 Bundeslaender + Gesamt row.
#  5) Sheet "Impfstoffdosen" (index 12): refreshed to a two-week comparison,
#     gains a new "Vorwoche" column (B), shifting the weekly figures to C.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Todesfälle und Fallsterblichkei (sheet 5) - Veraenderung column fixes
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("D2").Value  = "-8,5%"
$ws5.Range("D3").Value  = "-1,4%"
$ws5.Range("D4").Value  = "-5,2%"
$ws5.Range("D5").Value  = "-3,9%"
$ws5.Range("D7").Value  = "72,3%"
$ws5.Range("D8").Value  = "38,5%"
$ws5.Range("D9").Value  = "78,1%"
$ws5.Range("D10").Value = "59,4%"

# ---------------------------------------------------------------------------
# 2) Geimpfte Personen (sheet 9)
# ---------------------------------------------------------------------------
$ws9 = $wb.Worksheets.Item(9)
$ws9.Range("A1").Value = "Geimpfte Personen"
$ws9.Range("B1").Value = "Vorwoche"
$ws9.Range("C1").Value = "Stand 5.5."
$ws9.Range("D1").Value = "Anteil_Veraenderung"

$ws9.Range("A2").Value = "Gesamtbevölkerung"
$ws9.Range("B2").Value = ""
$ws9.Range("C2").Value = ""
$ws9.Range("D2").Value = ""

$ws9.Range("A3").Value = "Gesamt"
$ws9.Range("B3").Value = "20521685 (24,7 %)"
$ws9.Range("C3").Value = "24546919 (29,5 %)"
$ws9.Range("D3").Value = "4,8 PP"

$ws9.Range("A4").Value = "Nicht vollst. geimpft"
$ws9.Range("B4").Value = "14407957 (17,3 %)"
$ws9.Range("C4").Value = "17615335 (21,2 %)"
$ws9.Range("D4").Value = "3,9 PP"

$ws9.Range("A5").Value = "Vollst. geimpft"
$ws9.Range("B5").Value = "6113728 ( 7,4 %)"
$ws9.Range("C5").Value = "6931584 ( 8,3 %)"
$ws9.Range("D5").Value = "1,0 PP"

# ---------------------------------------------------------------------------
# 3) Impffortschritt (sheet 10)
# ---------------------------------------------------------------------------
$ws10 = $wb.Worksheets.Item(10)
$ws10.Range("A1").Value = "Impffortschritt"
$ws10.Range("B1").Value = "Vorwoche"
$ws10.Range("C1").Value = "letzteKW"
$ws10.Range("D1").Value = "Veraenderung"

$ws10.Range("A2").Value = "Impfungen pro Woche"
$ws10.Range("B2").Value = " "
$ws10.Range("C2").Value = " "
$ws10.Range("D2").Value = "   NA %"

$ws10.Range("A3").Value = "Gesamt"
$ws10.Range("B3").Value = "3493652"
$ws10.Range("C3").Value = "4637884"
$ws10.Range("D3").Value = " 32,8 %"

$ws10.Range("A4").Value = "davon in Impfzentren"
$ws10.Range("B4").Value = "2455305 ( 70,3 %)"
$ws10.Range("C4").Value = "2412430 ( 52 %)"
$ws10.Range("D4").Value = " -1,7 %"

$ws10.Range("A5").Value = "davon in ärztl. Praxen"
$ws10.Range("B5").Value = "1038347 ( 29,7 %)"
$ws10.Range("C5").Value = "2225454 ( 48 %)"
$ws10.Range("D5").Value = "114,3 %"

# ---------------------------------------------------------------------------
# 4) Regional Geimpfte (sheet 11) - refresh "Gesamt min. 1x" / "Gesamt vollst."
# ---------------------------------------------------------------------------
$ws11 = $wb.Worksheets.Item(11)

$regional = @(
    @{ Row = 2;  C = "29,5"; D = " 8,3" },
    @{ Row = 3;  C = "29,0"; D = " 8,0" },
    @{ Row = 4;  C = "30,0"; D = " 7,8" },
    @{ Row = 5;  C = "27,1"; D = "10,0" },
    @{ Row = 6;  C = "26,4"; D = " 8,0" },
    @{ Row = 7;  C = "30,5"; D = "10,4" },
    @{ Row = 8;  C = "29,7"; D = " 7,4" },
    @{ Row = 9;  C = "28,5"; D = " 8,0" },
    @{ Row = 10; C = "31,6"; D = " 7,4" },
    @{ Row = 11; C = "30,1"; D = " 7,4" },
    @{ Row = 12; C = "31,0"; D = " 7,9" },
    @{ Row = 13; C = "28,8"; D = " 9,4" },
    @{ Row = 14; C = "33,1"; D = " 9,0" },
    @{ Row = 15; C = "26,7"; D = "10,8" },
    @{ Row = 16; C = "29,6"; D = " 8,0" },
    @{ Row = 17; C = "27,8"; D = " 9,5" },
    @{ Row = 18; C = "27,4"; D = "12,2" }
)

foreach ($r in $regional) {
    $ws11.Range("C" + $r.Row).Value = $r.C
    $ws11.Range("D" + $r.Row).Value = $r.D
}

# ---------------------------------------------------------------------------
# 5) Impfstoffdosen (sheet 12)
# ---------------------------------------------------------------------------
$ws12 = $wb.Worksheets.Item(12)

$ws12.Range("A1").Value = "Impfstoffdosen"
$ws12.Range("B1").Value = "Vorwoche"
$ws12.Range("C1").Value = "dieseWoche"

$ws12.Range("A2").Value = "Biontech/Pfizer"
$ws12.Range("B2").Value = "19465427 (73,1 %)"
$ws12.Range("C2").Value = "23399097 (74,4 %)"

$ws12.Range("A3").Value = "Erstimpfungen"
$ws12.Range("B3").Value = "13756748"
$ws12.Range("C3").Value = "17014901"

$ws12.Range("A4").Value = "Zweitimpfungen"
$ws12.Range("B4").Value = "5708679"
$ws12.Range("C4").Value = "6384196"

$ws12.Range("A5").Value = "geliefert"
$ws12.Range("B5").Value = "23564774"
$ws12.Range("C5").Value = "25378274"

$ws12.Range("A6").Value = "Moderna"
$ws12.Range("B6").Value = "1507091 ( 5,7 %)"
$ws12.Range("C6").Value = "1932692 ( 6,1 %)"

$ws12.Range("A7").Value = "Erstimpfungen"
$ws12.Range("B7").Value = "1134048"
$ws12.Range("C7").Value = "1482621"

$ws12.Range("A8").Value = "Zweitimpfungen"
$ws12.Range("B8").Value = "373043"
$ws12.Range("C8").Value = "450071"

$ws12.Range("A9").Value = "geliefert"
$ws12.Range("B9").Value = "2742000"
$ws12.Range("C9").Value = "3118800"

$ws12.Range("A10").Value = "AstraZeneca"
$ws12.Range("B10").Value = "5661991 (21,3 %)"
$ws12.Range("C10").Value = "6127132 (19,5 %)"

$ws12.Range("A11").Value = "Erstimpfungen"
$ws12.Range("B11").Value = "5630437"
$ws12.Range("C11").Value = "6039606"

$ws12.Range("A12").Value = "Zweitimpfungen"
$ws12.Range("B12").Value = "31554"
$ws12.Range("C12").Value = "87526"

$ws12.Range("A13").Value = "geliefert"
$ws12.Range("B13").Value = "6899998"
$ws12.Range("C13").Value = "6959998"

$ws12.Range("A14").Value = "Johnson&Johnson"
$ws12.Range("B14").Value = "452 ( 0,0 %)"
$ws12.Range("C14").Value = "9791 ( 0,0 %)"

$ws12.Range("A15").Value = "geliefert"
$ws12.Range("B15").Value = "256800"
$ws12.Range("C15").Value = "256800"
